# Add a new "特殊类型" breakdown column to the tracking sheet.
#
# The existing "类型" column (C) is repurposed as "特殊类型", and a brand
# new "类型" column is inserted immediately to its right. The new column
# is populated with a rotating set of sub-categories (暑假单/专职单/
# 大学生单/特长单) for each data row. Everything that used to live in
# columns D (详细情况) and E (创建日期) simply shifts one column to the
# right, to E and F respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at D. Old D/E shift right to E/F; the new column
#    inherits column C's formatting (Excel's normal "insert column" rule).
$ws.Columns("D:D").Insert()

# 2) Relabel the old "类型" header (now still in C1) as "特殊类型", and
#    give the freshly inserted column its own "类型" header.
$ws.Range("C1").Value = "特殊类型"
$ws.Range("D1").Value = "类型"

# 3) Fill the new "类型" column (D2:D26) with the rotating sub-category
#    values, one per data row.
$subTypes = @("暑假单", "专职单", "大学生单", "特长单")
for ($row = 2; $row -le 26; $row++) {
    $ws.Cells.Item($row, 4).Value = $subTypes[($row - 2) % 4]
}

# 4) Match the new width of column D to column C, and leave the cursor
#    where the author's saved selection was.
$ws.Columns("D:D").ColumnWidth = $ws.Columns("C:C").ColumnWidth
$ws.Range("G2").Select()
